$wb = $excel.ActiveWorkbook

# Update the shared status text wherever "Ready for handoff" is used
# (Overview!B3, Overview!C3, zh-cn!C3, de-de!C3 all point at the same
# shared string, so a plain re-assignment of that text everywhere keeps
# them sharing one string).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("L3").Value = "Handback file name: j50ftwdg.rkk is different with handoff file name: 18db0254-b0ce-430f-8801-0743706e28a9.178cf6e5db8c6d9cd5e3a753cd8c900012bebb64.zh-cn."

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("L3").Value = "Handback file name: j50ftwdg.rkk is different with handoff file name: 18db0254-b0ce-430f-8801-0743706e28a9.178cf6e5db8c6d9cd5e3a753cd8c900012bebb64.de-de."
